# secard51/evidence.xlsx - "Add files via upload" commit replay
#
# Semantic edit: the long IBC-channel Juno address stored in the
# "JunoAddress" column (B2) of sheets A3 and A5 was trimmed down to the
# bare bech32 address (the "/channel-89/paloma" suffix was removed).
# Excel then re-saved the workbook, which also records wherever the
# cursor/selection was left on every sheet, and which sheet tab was
# active when the file was last saved.

$wb = $excel.ActiveWorkbook

$oldAddr = "juno1stv6sk0mvku34fj2mqrlyru6683866n306mfv52tlugtl322zmks26kg7a/channel-89/paloma"
$newAddr = "juno1stv6sk0mvku34fj2mqrlyru6683866n306mfv52tlugtl322zmks26kg7a"

# --- content edit: trim the Juno address on sheets A3 and A5 ---
$wsA3 = $wb.Worksheets.Item("A3")
$wsA3.Range("B2").Value = $newAddr

$wsA5 = $wb.Worksheets.Item("A5")
$wsA5.Range("B2").Value = $newAddr

# --- leftover cursor/selection state on each sheet ---
$wb.Worksheets.Item("A1").Range("A2").Select()
$wb.Worksheets.Item("A2").Range("A3").Select()
$wsA3.Range("A2").Select()
$wb.Worksheets.Item("A4").Range("B2").Select()
$wsA5.Range("B2").Select()
$wb.Worksheets.Item("A6").Range("B2").Select()
$wb.Worksheets.Item("A7").Range("A2").Select()
$wb.Worksheets.Item("A8").Range("A2").Select()

# --- active sheet moves from B5 to A9 (also restores its own selection) ---
$wsA9 = $wb.Worksheets.Item("A9")
$wsA9.Range("D9").Select()
$wsA9.Activate()
